{"js": "// Replace each arithmetic-problem answer cell's text with the new equation text.\n// The mapping below is old-equation-text -> new-equation-text, applied via\n// exact-text search + replace (each old string is unique within the document,\n// and none of the new strings collide with any other old string, so a single\n// left-to-right pass over the list is safe regardless of order).\nconst replacements = [\n  [\"71+8=79\", \"34+31=65\"],\n  [\"92-20=72\", \"12+62=74\"],\n  [\"42-20=22\", \"75+21=96\"],\n  [\"84-49=35\", \"42+37=79\"],\n  [\"53+2=55\", \"82-35=47\"],\n  [\"56-21=35\", \"60+39=99\"],\n  [\"74-61=13\", \"10-7=3\"],\n  [\"49-8=41\", \"14+76=90\"],\n  [\"76-74=2\", \"16+26=42\"],\n  [\"2+18=20\", \"39-15=24\"],\n  [\"93-45=48\", \"27-2=25\"],\n  [\"2+29=31\", \"46+35=81\"],\n  [\"14+55=69\", \"93-79=14\"],\n  [\"89-23=66\", \"17+11=28\"],\n  [\"42+6=48\", \"2+6=8\"],\n  [\"63+7=70\", \"5+13=18\"],\n  [\"86-9=77\", \"67+27=94\"],\n  [\"25+28=53\", \"97-34=63\"],\n  [\"29-22=7\", \"61-33=28\"],\n  [\"76-62=14\", \"98-92=6\"],\n  [\"40+7=47\", \"92-29=63\"],\n  [\"19+34=53\", \"52+32=84\"],\n  [\"90-54=36\", \"98-47=51\"],\n  [\"66+14=80\", \"26-25=1\"],\n  [\"69-64=5\", \"10+61=71\"],\n  [\"39+29=68\", \"48+11=59\"],\n  [\"12-2=10\", \"66-54=12\"],\n  [\"39+39=78\", \"9+14=23\"],\n  [\"46+5=51\", \"80-61=19\"],\n  [\"86-5=81\", \"90-3=87\"],\n  [\"29-5=24\", \"5+44=49\"],\n  [\"77-37=40\", \"43-4=39\"],\n  [\"63+27=90\", \"66-51=15\"],\n  [\"78-61=17\", \"9+11=20\"],\n  [\"60-2=58\", \"62-5=57\"],\n  [\"12+1=13\", \"16+19=35\"],\n  [\"86-45=41\", \"25+38=63\"],\n  [\"54+36=90\", \"79-19=60\"],\n  [\"85-4=81\", \"3+31=34\"],\n  [\"97-27=70\", \"35+10=45\"],\n  [\"35+28=63\", \"20+55=75\"],\n  [\"6+4=10\", \"70-62=8\"],\n  [\"30+9=39\", \"55+21=76\"],\n  [\"82-44=38\", \"38+27=65\"],\n  [\"37+57=94\", \"72+6=78\"],\n  [\"73-31=42\", \"62-38=24\"],\n  [\"35+20=55\", \"92-74=18\"],\n  [\"39-38=1\", \"64-8=56\"],\n  [\"46-5=41\", \"89-65=24\"],\n  [\"36-4=32\", \"8+70=78\"],\n  [\"13+69=82\", \"91-23=68\"],\n  [\"64+4=68\", \"98-30=68\"],\n  [\"10-9=1\", \"72-53=19\"],\n  [\"3+87=90\", \"5+48=53\"],\n  [\"68-47=21\", \"7+0=7\"],\n  [\"48-25=23\", \"51+38=89\"],\n  [\"64-35=29\", \"71-24=47\"],\n  [\"50-7=43\", \"23+23=46\"],\n  [\"83-43=40\", \"44+16=60\"],\n  [\"14+47=61\", \"25+48=73\"],\n  [\"29-18=11\", \"4+41=45\"],\n  [\"43+28=71\", \"19+14=33\"],\n  [\"44+53=97\", \"58-26=32\"],\n  [\"85-42=43\", \"72-22=50\"],\n  [\"74+23=97\", \"60-6=54\"],\n  [\"74-42=32\", \"50+49=99\"],\n  [\"19+41=60\", \"50-48=2\"],\n  [\"53+33=86\", \"18+76=94\"],\n  [\"74+24=98\", \"30-13=17\"],\n  [\"3+86=89\", \"81-72=9\"],\n  [\"18+1=19\", \"97-28=69\"],\n  [\"49-36=13\", \"33+4=37\"],\n  [\"73+22=95\", \"28+15=43\"],\n  [\"92-30=62\", \"35+7=42\"],\n  [\"64-56=8\", \"16+61=77\"],\n  [\"90-14=76\", \"21+46=67\"],\n  [\"3+95=98\", \"26-16=10\"],\n  [\"84-31=53\", \"66-52=14\"],\n  [\"50-3=47\", \"53+5=58\"],\n  [\"89-79=10\", \"45+18=63\"],\n  [\"91-43=48\", \"62-18=44\"],\n  [\"52-0=52\", \"96-95=1\"],\n  [\"55-19=36\", \"80-25=55\"],\n  [\"22+9=31\", \"11+8=19\"],\n  [\"75-29=46\", \"61+27=88\"],\n  [\"7+31=38\", \"14+3=17\"],\n  [\"83-2=81\", \"52-6=46\"],\n  [\"59+13=72\", \"47+29=76\"],\n  [\"49-26=23\", \"58-2=56\"],\n  [\"63-50=13\", \"14+5=19\"],\n  [\"55+38=93\", \"76-6=70\"],\n  [\"64-45=19\", \"28+44=72\"],\n  [\"66-1=65\", \"38+47=85\"],\n  [\"8+89=97\", \"93-89=4\"],\n  [\"33+27=60\", \"35+48=83\"],\n  [\"24-11=13\", \"79-6=73\"],\n  [\"68+14=82\", \"77-41=36\"],\n  [\"40+41=81\", \"62+20=82\"],\n  [\"83-9=74\", \"94-90=4\"],\n  [\"59-1=58\", \"16+55=71\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false, matchWildcards: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, 'Replace');\n  }\n}\nawait context.sync();", "ps1": "# Replace each arithmetic-problem answer cell's text with the new equation text.\n# The mapping below is old-equation-text -> new-equation-text, applied via\n# Word's Find/Replace (Content.Find.Execute) with MatchWildcards explicitly\n# disabled so '+' in the equations is treated literally. Every old string is\n# unique within the document and none of the new strings collide with any\n# other old string, so replacements can be applied in any order safely.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old='71+8=79'; New='34+31=65'},\n    @{Old='92-20=72'; New='12+62=74'},\n    @{Old='42-20=22'; New='75+21=96'},\n    @{Old='84-49=35'; New='42+37=79'},\n    @{Old='53+2=55'; New='82-35=47'},\n    @{Old='56-21=35'; New='60+39=99'},\n    @{Old='74-61=13'; New='10-7=3'},\n    @{Old='49-8=41'; New='14+76=90'},\n    @{Old='76-74=2'; New='16+26=42'},\n    @{Old='2+18=20'; New='39-15=24'},\n    @{Old='93-45=48'; New='27-2=25'},\n    @{Old='2+29=31'; New='46+35=81'},\n    @{Old='14+55=69'; New='93-79=14'},\n    @{Old='89-23=66'; New='17+11=28'},\n    @{Old='42+6=48'; New='2+6=8'},\n    @{Old='63+7=70'; New='5+13=18'},\n    @{Old='86-9=77'; New='67+27=94'},\n    @{Old='25+28=53'; New='97-34=63'},\n    @{Old='29-22=7'; New='61-33=28'},\n    @{Old='76-62=14'; New='98-92=6'},\n    @{Old='40+7=47'; New='92-29=63'},\n    @{Old='19+34=53'; New='52+32=84'},\n    @{Old='90-54=36'; New='98-47=51'},\n    @{Old='66+14=80'; New='26-25=1'},\n    @{Old='69-64=5'; New='10+61=71'},\n    @{Old='39+29=68'; New='48+11=59'},\n    @{Old='12-2=10'; New='66-54=12'},\n    @{Old='39+39=78'; New='9+14=23'},\n    @{Old='46+5=51'; New='80-61=19'},\n    @{Old='86-5=81'; New='90-3=87'},\n    @{Old='29-5=24'; New='5+44=49'},\n    @{Old='77-37=40'; New='43-4=39'},\n    @{Old='63+27=90'; New='66-51=15'},\n    @{Old='78-61=17'; New='9+11=20'},\n    @{Old='60-2=58'; New='62-5=57'},\n    @{Old='12+1=13'; New='16+19=35'},\n    @{Old='86-45=41'; New='25+38=63'},\n    @{Old='54+36=90'; New='79-19=60'},\n    @{Old='85-4=81'; New='3+31=34'},\n    @{Old='97-27=70'; New='35+10=45'},\n    @{Old='35+28=63'; New='20+55=75'},\n    @{Old='6+4=10'; New='70-62=8'},\n    @{Old='30+9=39'; New='55+21=76'},\n    @{Old='82-44=38'; New='38+27=65'},\n    @{Old='37+57=94'; New='72+6=78'},\n    @{Old='73-31=42'; New='62-38=24'},\n    @{Old='35+20=55'; New='92-74=18'},\n    @{Old='39-38=1'; New='64-8=56'},\n    @{Old='46-5=41'; New='89-65=24'},\n    @{Old='36-4=32'; New='8+70=78'},\n    @{Old='13+69=82'; New='91-23=68'},\n    @{Old='64+4=68'; New='98-30=68'},\n    @{Old='10-9=1'; New='72-53=19'},\n    @{Old='3+87=90'; New='5+48=53'},\n    @{Old='68-47=21'; New='7+0=7'},\n    @{Old='48-25=23'; New='51+38=89'},\n    @{Old='64-35=29'; New='71-24=47'},\n    @{Old='50-7=43'; New='23+23=46'},\n    @{Old='83-43=40'; New='44+16=60'},\n    @{Old='14+47=61'; New='25+48=73'},\n    @{Old='29-18=11'; New='4+41=45'},\n    @{Old='43+28=71'; New='19+14=33'},\n    @{Old='44+53=97'; New='58-26=32'},\n    @{Old='85-42=43'; New='72-22=50'},\n    @{Old='74+23=97'; New='60-6=54'},\n    @{Old='74-42=32'; New='50+49=99'},\n    @{Old='19+41=60'; New='50-48=2'},\n    @{Old='53+33=86'; New='18+76=94'},\n    @{Old='74+24=98'; New='30-13=17'},\n    @{Old='3+86=89'; New='81-72=9'},\n    @{Old='18+1=19'; New='97-28=69'},\n    @{Old='49-36=13'; New='33+4=37'},\n    @{Old='73+22=95'; New='28+15=43'},\n    @{Old='92-30=62'; New='35+7=42'},\n    @{Old='64-56=8'; New='16+61=77'},\n    @{Old='90-14=76'; New='21+46=67'},\n    @{Old='3+95=98'; New='26-16=10'},\n    @{Old='84-31=53'; New='66-52=14'},\n    @{Old='50-3=47'; New='53+5=58'},\n    @{Old='89-79=10'; New='45+18=63'},\n    @{Old='91-43=48'; New='62-18=44'},\n    @{Old='52-0=52'; New='96-95=1'},\n    @{Old='55-19=36'; New='80-25=55'},\n    @{Old='22+9=31'; New='11+8=19'},\n    @{Old='75-29=46'; New='61+27=88'},\n    @{Old='7+31=38'; New='14+3=17'},\n    @{Old='83-2=81'; New='52-6=46'},\n    @{Old='59+13=72'; New='47+29=76'},\n    @{Old='49-26=23'; New='58-2=56'},\n    @{Old='63-50=13'; New='14+5=19'},\n    @{Old='55+38=93'; New='76-6=70'},\n    @{Old='64-45=19'; New='28+44=72'},\n    @{Old='66-1=65'; New='38+47=85'},\n    @{Old='8+89=97'; New='93-89=4'},\n    @{Old='33+27=60'; New='35+48=83'},\n    @{Old='24-11=13'; New='79-6=73'},\n    @{Old='68+14=82'; New='77-41=36'},\n    @{Old='40+41=81'; New='62+20=82'},\n    @{Old='83-9=74'; New='94-90=4'},\n    @{Old='59-1=58'; New='16+55=71'},\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $result = $find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"Find/Replace failed for: $($pair.Old)\"\n    }\n}\n"}
